$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 26
$ws.Range("I2").Value = 85
$ws.Range("J2").Value = 375
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 108
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 57
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = 36
$ws.Range("T2").Value = 67
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 596
$ws.Range("X2").Value = 570
$ws.Range("Z2").Value = 8
$ws.Range("AA2").Value = 8
